$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-18 09:47:08", 0.001),
    @("2023-12-18 09:48:06", 0.0044),
    @("2023-12-18 09:49:17", 0.006),
    @("2023-12-18 09:49:22", 0.0004),
    @("2023-12-18 09:49:32", 0.0004)
)

$startRow = 424
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
